$d = $word.ActiveDocument

# --- Paragraph 1: split "This is a Microsoft word document." into
#     the original text (plus two trailing spaces) followed by three
#     red-colored runs: "(This is a change " [en dash] "Ve",
#     "rsion for main branch", ")"
$p1 = $d.Paragraphs(1).Range
$pEnd = $p1.End - 1   # position just before the paragraph mark

$ins = $d.Range($pEnd, $pEnd)
$ins.InsertAfter("  ")
$pEnd = $pEnd + 2

$ins = $d.Range($pEnd, $pEnd)
$ins.InsertAfter([string]::Concat("(This is a change ", [char]0x2013, " Ve"))
$r2start = $pEnd
$pEnd = $pEnd + 22
$r2 = $d.Range($r2start, $pEnd)
$r2.Font.Color = 255

$ins = $d.Range($pEnd, $pEnd)
$ins.InsertAfter("rsion for main branch")
$r3start = $pEnd
$pEnd = $pEnd + 21
$r3 = $d.Range($r3start, $pEnd)
$r3.Font.Color = 255

$ins = $d.Range($pEnd, $pEnd)
$ins.InsertAfter(")")
$r4start = $pEnd
$pEnd = $pEnd + 1
$r4 = $d.Range($r4start, $pEnd)
$r4.Font.Color = 255

Write-Output $d.Paragraphs(1).Range.Text
